$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bmwModels")

# Turn the single "Models" column into a two-column Name/Model table:
#  - A1 header "Models" -> "Name"; new header B1 "Model" (bold, matching A1)
#  - A2 "X1" -> "BMW"; the X1 value moves over to new column B2
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Model"
$ws.Range("B1").Font.Bold = $true
$ws.Range("A2").Value = "BMW"
$ws.Range("B2").Value = "X1"

# Switch focus to this sheet and select A3, as captured in the saved view state
$ws.Activate() | Out-Null
$ws.Range("A3").Select() | Out-Null
